# Auto-generated Excel COM-interop edit script
# Applies updated horarios (schedule) data per commit diff

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")

# --- Header updates (timestamp / total rows) and cell swaps/edits ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 20:33:25"
$ws1.Cells.Item(3,1).Value = "Total filas: 545"
$ws1.Cells.Item(54,1).Value = "06:02:16"
$ws1.Cells.Item(54,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(54,4).Value = 74
$ws1.Cells.Item(55,1).Value = "06:37:24"
$ws1.Cells.Item(55,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(55,4).Value = 39
$ws1.Cells.Item(140,1).Value = "09:25:30"
$ws1.Cells.Item(140,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(140,4).Value = 58
$ws1.Cells.Item(141,1).Value = "10:11:11"
$ws1.Cells.Item(141,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(141,4).Value = 12
$ws1.Cells.Item(209,3).Value = "17_ROMERO"
$ws1.Cells.Item(210,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(225,1).Value = "11:15:53"
$ws1.Cells.Item(225,3).Value = "17_ROMERO"
$ws1.Cells.Item(225,4).Value = 91
$ws1.Cells.Item(226,1).Value = "12:24:14"
$ws1.Cells.Item(226,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(226,4).Value = 22
$ws1.Cells.Item(259,1).Value = "12:57:33"
$ws1.Cells.Item(259,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(259,4).Value = 53
$ws1.Cells.Item(260,1).Value = "11:58:46"
$ws1.Cells.Item(260,3).Value = "215A_EL PATO"
$ws1.Cells.Item(260,4).Value = 112
$ws1.Cells.Item(295,1).Value = "14:17:27"
$ws1.Cells.Item(295,3).Value = "10_OLMOS"
$ws1.Cells.Item(295,4).Value = 43
$ws1.Cells.Item(296,1).Value = "13:24:27"
$ws1.Cells.Item(296,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(296,4).Value = 96
$ws1.Cells.Item(328,3).Value = "17_ROMERO"
$ws1.Cells.Item(330,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(338,1).Value = "14:42:26"
$ws1.Cells.Item(338,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(338,4).Value = 83
$ws1.Cells.Item(339,1).Value = "15:53:26"
$ws1.Cells.Item(339,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(339,4).Value = 12
$ws1.Cells.Item(346,1).Value = "16:14:21"
$ws1.Cells.Item(346,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(346,4).Value = 5
$ws1.Cells.Item(347,1).Value = "15:21:47"
$ws1.Cells.Item(347,3).Value = "215C_EL PATO"
$ws1.Cells.Item(347,4).Value = 58
$ws1.Cells.Item(362,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(363,1).Value = "16:14:21"
$ws1.Cells.Item(363,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(363,4).Value = 29
$ws1.Cells.Item(364,1).Value = "14:56:04"
$ws1.Cells.Item(364,3).Value = "225_GOMEZ"
$ws1.Cells.Item(364,4).Value = 107
$ws1.Cells.Item(377,1).Value = "16:14:21"
$ws1.Cells.Item(377,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(377,4).Value = 51
$ws1.Cells.Item(378,1).Value = "16:47:11"
$ws1.Cells.Item(378,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(378,4).Value = 18
$ws1.Cells.Item(379,3).Value = "215A_EL PATO"
$ws1.Cells.Item(380,1).Value = "16:30:20"
$ws1.Cells.Item(380,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(380,4).Value = 37
$ws1.Cells.Item(381,1).Value = "16:39:47"
$ws1.Cells.Item(381,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(381,4).Value = 28
$ws1.Cells.Item(499,3).Value = "14_ABASTO"
$ws1.Cells.Item(500,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(516,1).Value = "19:13:07"
$ws1.Cells.Item(516,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(516,4).Value = 69
$ws1.Cells.Item(517,1).Value = "18:34:43"
$ws1.Cells.Item(517,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(517,4).Value = 108
$ws1.Cells.Item(522,1).Value = "20:33:25"
$ws1.Cells.Item(522,2).Value = "20:33"
$ws1.Cells.Item(522,3).Value = "17_ROMERO"
$ws1.Cells.Item(522,4).Value = 0
$ws1.Cells.Item(523,1).Value = "20:33:25"
$ws1.Cells.Item(523,2).Value = "20:34"
$ws1.Cells.Item(523,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(523,4).Value = 1
$ws1.Cells.Item(524,1).Value = "19:13:07"
$ws1.Cells.Item(524,2).Value = "20:44"
$ws1.Cells.Item(524,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(524,4).Value = 91
$ws1.Cells.Item(525,1).Value = "20:33:25"
$ws1.Cells.Item(525,2).Value = "20:46"
$ws1.Cells.Item(525,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(525,4).Value = 13
$ws1.Cells.Item(526,2).Value = "20:52"
$ws1.Cells.Item(526,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(526,4).Value = 74
$ws1.Cells.Item(527,1).Value = "19:38:38"
$ws1.Cells.Item(527,2).Value = "20:52"
$ws1.Cells.Item(527,3).Value = "15_ABASTO"
$ws1.Cells.Item(527,4).Value = 74
$ws1.Cells.Item(528,1).Value = "19:51:02"
$ws1.Cells.Item(528,2).Value = "20:53"
$ws1.Cells.Item(528,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(528,4).Value = 62
$ws1.Cells.Item(529,1).Value = "19:38:38"
$ws1.Cells.Item(529,2).Value = "20:56"
$ws1.Cells.Item(529,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(529,4).Value = 78
$ws1.Cells.Item(530,1).Value = "19:13:07"
$ws1.Cells.Item(530,2).Value = "20:56"
$ws1.Cells.Item(530,3).Value = "10_OLMOS"
$ws1.Cells.Item(530,4).Value = 103
$ws1.Cells.Item(531,2).Value = "20:57"
$ws1.Cells.Item(531,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(531,4).Value = 104
$ws1.Cells.Item(532,1).Value = "19:13:07"
$ws1.Cells.Item(532,2).Value = "21:04"
$ws1.Cells.Item(532,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(532,4).Value = 111
$ws1.Cells.Item(533,1).Value = "20:33:25"
$ws1.Cells.Item(533,2).Value = "21:04"
$ws1.Cells.Item(533,3).Value = "15_ABASTO"
$ws1.Cells.Item(533,4).Value = 31
$ws1.Cells.Item(534,2).Value = "21:07"
$ws1.Cells.Item(534,3).Value = "215B_EL PATO"
$ws1.Cells.Item(534,4).Value = 89
$ws1.Cells.Item(535,1).Value = "19:13:07"
$ws1.Cells.Item(535,2).Value = "21:08"
$ws1.Cells.Item(535,3).Value = "215B_EL PATO"
$ws1.Cells.Item(535,4).Value = 115
$ws1.Cells.Item(536,1).Value = "19:38:38"
$ws1.Cells.Item(536,2).Value = "21:20"
$ws1.Cells.Item(536,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(536,4).Value = 102
$ws1.Cells.Item(537,2).Value = "21:21"
$ws1.Cells.Item(537,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(537,4).Value = 90
$ws1.Cells.Item(538,1).Value = "19:38:38"
$ws1.Cells.Item(538,2).Value = "21:22"
$ws1.Cells.Item(538,3).Value = "10_OLMOS"
$ws1.Cells.Item(538,4).Value = 104
$ws1.Cells.Item(539,2).Value = "21:23"
$ws1.Cells.Item(539,3).Value = "15_ABASTO"
$ws1.Cells.Item(539,4).Value = 71
$ws1.Cells.Item(540,2).Value = "21:23"
$ws1.Cells.Item(540,3).Value = "10_OLMOS"
$ws1.Cells.Item(540,4).Value = 92

# --- New rows appended at the end of LP1912 ---
$ws1.Cells.Item(541,1).Value = "20:33:25"
$ws1.Cells.Item(541,2).Value = "21:32"
$ws1.Cells.Item(541,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(541,4).Value = 59
$ws1.Cells.Item(541,5).Value = "LP1912"
$ws1.Cells.Item(542,1).Value = "19:51:02"
$ws1.Cells.Item(542,2).Value = "21:38"
$ws1.Cells.Item(542,3).Value = "17_ROMERO"
$ws1.Cells.Item(542,4).Value = 107
$ws1.Cells.Item(542,5).Value = "LP1912"
$ws1.Cells.Item(543,1).Value = "19:51:02"
$ws1.Cells.Item(543,2).Value = "21:38"
$ws1.Cells.Item(543,3).Value = "14_ABASTO"
$ws1.Cells.Item(543,4).Value = 107
$ws1.Cells.Item(543,5).Value = "LP1912"
$ws1.Cells.Item(544,1).Value = "20:12:07"
$ws1.Cells.Item(544,2).Value = "21:40"
$ws1.Cells.Item(544,3).Value = "17_ROMERO"
$ws1.Cells.Item(544,4).Value = 88
$ws1.Cells.Item(544,5).Value = "LP1912"
$ws1.Cells.Item(545,1).Value = "19:51:02"
$ws1.Cells.Item(545,2).Value = "21:47"
$ws1.Cells.Item(545,3).Value = "215A_EL PATO"
$ws1.Cells.Item(545,4).Value = 116
$ws1.Cells.Item(545,5).Value = "LP1912"
$ws1.Cells.Item(546,1).Value = "20:12:07"
$ws1.Cells.Item(546,2).Value = "21:53"
$ws1.Cells.Item(546,3).Value = "10_OLMOS"
$ws1.Cells.Item(546,4).Value = 101
$ws1.Cells.Item(546,5).Value = "LP1912"
$ws1.Cells.Item(547,1).Value = "20:33:25"
$ws1.Cells.Item(547,2).Value = "22:01"
$ws1.Cells.Item(547,3).Value = "17_ROMERO"
$ws1.Cells.Item(547,4).Value = 88
$ws1.Cells.Item(547,5).Value = "LP1912"
$ws1.Cells.Item(548,1).Value = "20:33:25"
$ws1.Cells.Item(548,2).Value = "22:08"
$ws1.Cells.Item(548,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(548,4).Value = 95
$ws1.Cells.Item(548,5).Value = "LP1912"
$ws1.Cells.Item(549,1).Value = "20:33:25"
$ws1.Cells.Item(549,2).Value = "22:19"
$ws1.Cells.Item(549,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(549,4).Value = 106
$ws1.Cells.Item(549,5).Value = "LP1912"
$ws1.Cells.Item(550,1).Value = "20:33:25"
$ws1.Cells.Item(550,2).Value = "22:28"
$ws1.Cells.Item(550,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(550,4).Value = 115
$ws1.Cells.Item(550,5).Value = "LP1912"

# --- LP1912-215 sheet: timestamp-only update ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 20:33:25"

# --- 6203-6173 sheet: timestamp/total update + new row ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 20:33:25"
$ws3.Cells.Item(3,1).Value = "Total filas: 67"
$ws3.Cells.Item(72,1).Value = "20:33:25"
$ws3.Cells.Item(72,2).Value = "22:20"
$ws3.Cells.Item(72,3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(72,4).Value = 107
$ws3.Cells.Item(72,5).Value = "L6173"
